$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32, shifting existing rows (32..119) down to (33..120).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new data record.
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44607
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112001
$ws.Range("G32").Value = "Berenjena"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 480
$ws.Range("K32").Value = 8500
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = 8750
$ws.Range("N32").Value = '$/caja 50 unidades'
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 175
$ws.Range("Q32").Value = 50
$ws.Range("R32").Value = "Hortaliza"
